{"js": "// Remove decorative inline-picture paragraphs and the leftover empty\n// spacer paragraphs (w:spacing w:before=\"40\") that used to sit between a\n// code-block table and the paragraph following it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load what we need to decide which paragraphs must go: text, the\n// \"space before\" paragraph formatting (in points - 2pt == 40 twips) and\n// whether the paragraph hosts an inline picture.\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text,spaceBefore\");\n  p.inlinePictures.load(\"items\");\n}\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const hasPicture = p.inlinePictures.items.length > 0;\n  const isEmpty = p.text === \"\";\n\n  if (hasPicture) {\n    // Paragraph whose only content is an inline image.\n    toDelete.push(p);\n  } else if (isEmpty && p.spaceBefore === 2) {\n    // Empty spacer paragraph (<w:spacing w:before=\"40\"/>) directly after\n    // a table - 40 twentieths of a point == 2pt.\n    toDelete.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove decorative inline-picture paragraphs and the leftover empty\n# spacer paragraphs (w:spacing w:before=\"40\", 2pt) that used to sit\n# between a code-block table and the paragraph following it.\n\n$d = $word.ActiveDocument\n$n = $d.Paragraphs.Count\n\n# First pass: identify the Start/End of every paragraph that must be\n# removed. We collect plain numbers (not live Range objects) because the\n# document positions shift only *after* we actually delete something, and\n# we will delete from the end of the document backwards so earlier,\n# already-recorded offsets stay valid.\n$toDelete = @()\nfor ($i = 1; $i -le $n; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n\n    $hasPicture = $r.InlineShapes.Count -gt 0\n    # An untouched paragraph mark reads back as just \"`r\"; the phantom\n    # zero-length paragraph Word reports for each table's own mark is\n    # excluded automatically because its SpaceBefore differs / it is not\n    # a real empty-before-table spacer (guarded by the Start<End check).\n    $isEmptySpacer = (-not $hasPicture) -and ($r.Text -eq \"`r\") -and ($p.SpaceBefore -eq 2) -and ($r.End -gt $r.Start)\n\n    if ($hasPicture -or $isEmptySpacer) {\n        $toDelete += , @($r.Start, $r.End)\n    }\n}\n\n# Second pass: delete in reverse document order so already-collected\n# offsets for earlier paragraphs remain valid.\nfor ($k = $toDelete.Count - 1; $k -ge 0; $k--) {\n    $pair = $toDelete[$k]\n    $rng = $d.Range($pair[0], $pair[1])\n    $rng.Delete()\n}\n"}
